$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New incident rows (176-182) appended to the log, all stored as plain text
# to match the existing sheet's inline-string formatting (dates/times kept
# as text, not converted to Excel date/time serials).
$rows = @(
    @{ r = 176; A = "WC49 P5H";     B = "Screw K30 no lo detecta puesto"; C = "2024-06-10"; D = "11:47:59"; E = "Mañana"; F = "11:48:00"; G = "0:00:01"; H = "-0.00 minutos" }
    @{ r = 177; A = "WC49 P5H";     B = "No lee QR";                      C = "2024-06-10"; D = "11:48:06"; E = "Mañana"; F = "11:48:07"; G = "0:00:01"; H = "0.03 minutos" }
    @{ r = 178; A = "WC49 P5H";     B = "Atasco tuerca";                  C = "2024-06-10"; D = "11:48:44"; E = "Mañana"; F = "11:48:46"; G = "0:00:02"; H = "0.12 minutos" }
    @{ r = 179; A = "WV50 FILTER";  B = "QR desplazado";                  C = "2024-06-10"; D = "11:51:14"; E = "Mañana"; F = "11:51:15"; G = "0:00:01"; H = "-0.01 minutos" }
    @{ r = 180; A = "WC48 P5F";     B = "AOI (fallo etiqueta)";           C = "2024-06-10"; D = "11:57:17"; E = "Mañana"; F = "11:57:17"; G = "0:00:00"; H = "-0.00 minutos" }
    @{ r = 181; A = "WC48 P5F";     B = "Etiquetadora";                   C = "2024-06-10"; D = "11:57:29"; E = "Mañana"; F = "11:57:29"; G = "0:00:00"; H = "0.05 minutos" }
    @{ r = 182; A = "WC47 NACP";    B = "No atornilla tapa";              C = "2024-06-10"; D = "12:01:55"; E = "Mañana"; F = $null;       G = $null;       H = "-0.00 minutos" }
)

foreach ($row in $rows) {
    $r = $row.r

    # Force text format on the whole row first so values (dates, times,
    # durations) are preserved verbatim as strings, not auto-converted.
    $ws.Range("A$r`:H$r").NumberFormat = "@"

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    if ($null -ne $row.F) { $ws.Cells.Item($r, 6).Value = $row.F }
    if ($null -ne $row.G) { $ws.Cells.Item($r, 7).Value = $row.G }
    $ws.Cells.Item($r, 8).Value = $row.H
}
